$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 vote-count columns (H through AA) with the refreshed election results
$ws.Range("H2").Value  = 98
$ws.Range("I2").Value  = 298
$ws.Range("J2").Value  = 1216
$ws.Range("K2").Value  = 5
$ws.Range("L2").Value  = 310
$ws.Range("M2").Value  = 17
$ws.Range("N2").Value  = 227
$ws.Range("O2").Value  = 2
$ws.Range("P2").Value  = 7
$ws.Range("Q2").Value  = 3
$ws.Range("R2").Value  = 12
$ws.Range("S2").Value  = 125
$ws.Range("T2").Value  = 206
$ws.Range("U2").Value  = 12
$ws.Range("V2").Value  = 1899
$ws.Range("W2").Value  = 0
$ws.Range("X2").Value  = 1789
$ws.Range("Y2").Value  = 3
$ws.Range("Z2").Value  = 23
$ws.Range("AA2").Value = 15
